$d = $word.ActiveDocument

# The three distinct body paragraphs that get repeated through the new block.
$t1 = 'Jury mast grapple strike colors cutlass fore spike loaded to the gunwalls Gold Road Jolly Roger skysail. Avast Jack Tar hail-shot scuttle spanker reef parley aye swab gangway. Lee black jack holystone Admiral of the Black mutiny mizzenmast Nelsons folly swing the lead warp bucko.'
$t2 = 'Smartly crimp squiffy Sail ho schooner Sink me lass overhaul pressgang piracy. Brethren of the Coast trysail marooned run a shot across the bow cog coffer Pirate Round provost yo-ho-ho grog blossom. Brigantine bring a spring upon her cable Privateer parrel marooned booty tackle heave down Yellow Jack pressgang.'
$t3 = 'Keel cog long clothes no prey, no pay lugsail six pounders Shiver me timbers run a shot across the bow belaying pin plunder. Port loot Plate Fleet scourge of the seven seas dance the hempen jig gabion Jack Tar Cat o''nine tails wherry cutlass. Lugger Corsair red ensign careen deadlights chase warp pillage Sink me Jolly Roger.'

# Sequence of 32 new paragraphs to insert right after the "Trysail Sail ho ..."
# paragraph: (Empty, t1, Empty, t2, Empty, t3) repeated 5 times, then 2 more
# Empty paragraphs, mirroring the rest of the document's blank-line-between-
# paragraphs pattern.
$sequence = @(
    $null, $t1, $null, $t2, $null, $t3,
    $null, $t1, $null, $t2, $null, $t3,
    $null, $t1, $null, $t2, $null, $t3,
    $null, $t1, $null, $t2, $null, $t3,
    $null, $t1, $null, $t2, $null, $t3,
    $null, $null
)

# Find the "Trysail Sail ho ..." paragraph - the anchor point for the insert.
$anchor = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "Trysail Sail ho*") {
        $anchor = $p
    }
}
if ($anchor -eq $null) {
    $anchor = $d.Paragraphs.Item(9)
}

$cur = $anchor

foreach ($txt in $sequence) {
    $r = $cur.Range
    $r.Collapse(0)
    $r.InsertParagraphAfter()
    $cur = $cur.Next()
    if ($txt) {
        $cur.Range.InsertAfter($txt)
    }
}

Write-Host $d.Paragraphs.Count
